# Rename data elements to underscore_case and add the missing
# "participant_response" element (which duplicates key_press).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) pairTrialsPerBlock -> pair_trials_per_block (row 13, column A)
$ws.Cells.Item(13, 1).Value = "pair_trials_per_block"

# 2) Insert a new data-definition row before the old "accuracy" row (30),
#    pushing it down to row 31, and fill in the new row describing
#    "participant_response".
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 1).Value = "participant_response"
$ws.Cells.Item(30, 2).Value = "Int"
$ws.Cells.Item(30, 3).Value = "participant response, same as key_press"
$ws.Cells.Item(30, 4).Value = "ascii value"
$ws.Cells.Item(30, 5).Value = "'49,50"

# Resize the ElementName / ElementDescription columns to fit the new,
# longer content.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Leave the selection where the edit finished.
$ws.Range("D30").Select() | Out-Null
